$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Packing list")

# Insert 5 new rows starting at row 20 (pushes old row 20 down to row 25,
# and the previous content that began there cascades further down).
$ws.Rows("20:24").Insert()

# Match row 19's height on the freshly inserted rows
$ws.Rows("20:24").RowHeight = $ws.Rows("19").RowHeight

# The printable area grew by the same 5 rows (was $A$1:$I$30)
$ws.PageSetup.PrintArea = '$A$1:$I$35'

# Restore the window/selection to where the editor ended up working
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 8
$ws.Range("B21").Select() | Out-Null

